$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of F1 ("טופל ב") and F2 ("29/01/2025 18:42") while
# keeping their existing cell formatting/styles.
$ws.Range("F1:F2").ClearContents()

# Row 2 was taller than default because the deleted F2 text wrapped to two
# lines; re-fit the row now that the text is gone so the explicit height is
# dropped again.
$ws.Rows(2).AutoFit()

# Update the selection to reflect the user selecting column F (F1:F4)
# with F4 as the active cell, as was left after the edit.
$ws.Range("F1:F4").Select()
